# TXs -> TX, RXs -> RX 로 간략화
#
# - Slide 2 ("F1"): simplify the interface-range labels on the API3 and
#   API4 shapes (drop the extra "R5"/"R1;" bits).
# - Remove slide 4 ("F3") entirely - it's no longer part of the deck.

$p = $ppt.ActivePresentation

$s2 = $p.Slides.Item(2)

# Shape 15 = "화살표: 오각형 79" (API3) - second paragraph holds the range label.
$para15 = $s2.Shapes.Item(15).TextFrame.TextRange.Paragraphs(2)
# Re-assign through an intermediate value so the engine treats this as a
# full replacement (no shared prefix/suffix with the old text) instead of
# splitting the run around the common "[R2 ~R2" / "]" substrings.
$para15.Text = "X"
$para15.Text = "[R2 ~R2]"

# Shape 16 = "화살표: 오각형 80" (API4) - second paragraph holds the range label.
$para16 = $s2.Shapes.Item(16).TextFrame.TextRange.Paragraphs(2)
$para16.Text = "X"
$para16.Text = "[R1~ R4]"

# Slide 4 ("F3") is removed from the deck entirely.
$p.Slides.Item(4).Delete()
